# Add a new "2022-Q3" sheet before the current second sheet ("2022-Q2"),
# fill it with the quarterly fund-holding table, and insert a matching
# summary row at the top of the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet in the right tab position.
#    Grab formatting references from the existing "2022-Q2" sheet
#    *before* anything is renamed/reordered, so the new sheet reuses
#    the same style indexes (bold+border header / index column) instead
#    of minting new ones.
# ---------------------------------------------------------------------
$fmtSrc = $wb.Worksheets.Item("2022-Q2")
$target = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($target)
$newSheet.Name = "2022-Q3"

$fmtSrc.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$fmtSrc.Range("A2").Copy()
$newSheet.Range("A2:A44").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Populate the new sheet's header row text.
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 3. Populate the data rows (A=index, B=code, C=name, D..G=text figures,
#    H=rank number). Row 44 col G is a genuine number (0), matching the
#    source data; every other D/E/F/G value is kept as text.
# ---------------------------------------------------------------------
$rows = @(
  @("515220","国泰中证煤炭ETF","47.41","99.43","7.97","3.7786",4),
  @("510500","南方中证500ETF","387.70","95.58","0.56","2.1711",5),
  @("161032","富国中证煤炭指数A","18.28","94.24","7.48","1.3673",4),
  @("168204","中融中证煤炭指数A","9.23","92.26","7.29","0.6729",4),
  @("013275","富国中证煤炭指数C","8.56","94.24","7.48","0.6403",4),
  @("159922","嘉实中证500ETF","64.35","98.72","0.59","0.3797",5),
  @("159930","汇添富中证能源ETF","4.30","99.45","5.75","0.2472",7),
  @("512500","华夏中证500ETF","36.51","98.13","0.58","0.2118",5),
  @("510510","广发中证500ETF","23.97","97.86","0.58","0.1390",5),
  @("510410","博时上证自然资源ETF","4.47","98.01","2.94","0.1314",10),
  @("159820","天弘中证500ETF","22.24","95.49","0.57","0.1268",5),
  @("510580","易方达中证500ETF","18.94","94.41","0.56","0.1061",5),
  @("620001","金元顺安宝石动力混合","1.01","40.12","5.93","0.0599",2),
  @("159982","鹏华中证500ETF","10.21","95.57","0.57","0.0582",5),
  @("159968","博时中证500ETF","7.26","94.95","0.56","0.0407",5),
  @("002316","创金合信中证500指数增强C","2.51","93.93","1.35","0.0339",1),
  @("002311","创金合信中证500指数增强A","2.28","93.93","1.35","0.0308",1),
  @("510590","平安中证500ETF","5.30","96.44","0.58","0.0307",5),
  @("512510","华泰柏瑞中证500ETF","4.91","95.72","0.56","0.0275",5),
  @("160616","鹏华中证500指数（LOF）A","4.94","92.84","0.55","0.0272",5),
  @("159945","广发中证全指能源ETF","0.52","97.88","4.36","0.0227",7),
  @("510530","工银中证500ETF","2.58","98.55","0.59","0.0152",5),
  @("510560","国寿安保中证500ETF","2.23","99.38","0.59","0.0132",5),
  @("165511","信诚中证500指数（LOF）A","2.23","92.39","0.56","0.0125",5),
  @("561350","国泰中证500ETF","1.85","98.49","0.57","0.0105",5),
  @("010992","西藏东财中证500指数A","1.58","94.55","0.56","0.0088",5),
  @("510440","大成中证500沪市ETF","0.40","97.01","1.10","0.0044",5),
  @("515190","中银证券中证500ETF","0.64","99.17","0.59","0.0038",5),
  @("515530","泰康中证500ETF","0.63","96.14","0.57","0.0036",5),
  @("010993","西藏东财中证500指数C","0.58","94.55","0.56","0.0032",5),
  @("159935","景顺长城中证500ETF","0.53","98.38","0.58","0.0031",5),
  @("660011","农银中证500指数","0.54","94.17","0.56","0.0030",5),
  @("008113","中泰中证500指数增强C","0.44","91.57","0.55","0.0024",10),
  @("159999","永赢中证500ETF","0.42","95.65","0.58","0.0024",5),
  @("006611","人保中证500指数","0.39","92.75","0.58","0.0023",9),
  @("006938","鹏华中证500指数（LOF）C","0.40","92.84","0.55","0.0022",5),
  @("515550","中融中证500ETF","0.33","93.85","0.57","0.0019",5),
  @("519117","浦银安盛基本面400指数","0.22","91.78","0.65","0.0014",3),
  @("510550","方正富邦中证500ETF","0.19","97.10","0.58","0.0011",6),
  @("008112","中泰中证500指数增强A","0.16","91.57","0.55","0.0009",10),
  @("013119","信诚中证500指数（LOF）C","0.10","92.39","0.56","0.0006",5),
  @("510570","兴业中证500ETF","0.09","95.85","0.60","0.0005",8),
  @("016814","中融中证煤炭指数C","0.00","92.26","7.29",0,4)
)

# Text-format columns B-G up front so numeric-looking strings (fund codes
# with leading zeros, "47.41", "0.0024", ...) are stored as text, not
# auto-converted to numbers.
$lastRow = $rows.Length + 1
$newSheet.Range("B2:G$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# Row 44's G value (index 42) is a genuine number in the source file, so
# restore its native numeric format/type after the blanket text pass.
$newSheet.Range("G$lastRow").NumberFormat = "General"
$newSheet.Range("G$lastRow").Value = 0

# B-G data cells have no explicit style in the source workbook; drop the
# "@" text format back off now that the values are safely stored as text.
$newSheet.Range("B2:G$lastRow").Style = "Normal"

# ---------------------------------------------------------------------
# 4. Insert the matching summary row at the top of the "总计" table.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows(2).Insert()
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 43
$summary.Cells.Item(2, 4).Value = 10.4
$summary.Range("B2:D2").Style = "Normal"

$summary.Cells.Item(2, 1).Value = 0
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
